# Removing the date column for now. The Excel type provider reports an
# error (relates to issue #236). Drop the "Date" header/value column
# (D1:D2) from Sheet1, leaving D2's existing number-format style in place
# but with no formula/value, and rename the default cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# D1 held the "Date" header (shared string) - remove it outright.
$ws.Range("D1").ClearContents()

# D2 held "=NOW()" formatted with the date-time style (s="3"). Clear the
# formula/value but keep the cell (and its style) in place.
$ws.Range("D2").ClearContents()

# Reflect the column having been the user's last point of interest.
$ws.Range("D1:D2").Select()

# Rename the built-in "Normal" cell style to "Standard".
$style = $wb.Styles.Item("Normal")
$style.Name = "Standard"
